$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily data row needs to be inserted at row 701 (for 2026/01/22, 06:00),
# pushing all the subsequent rows (old 701..742) down by one (to 702..743).
$ws.Rows.Item(701).Insert()

# Make sure column A keeps being stored as plain text (the sheet stores dates
# as literal strings like "2026/01/22", not as real Excel date serials).
$ws.Range("A701").NumberFormat = "@"

$ws.Range("A701").Value = "2026/01/22"
$ws.Range("B701").Value = "木"
$ws.Range("C701").Value = 6
$ws.Range("D701").Value = 163
